$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'medical compression shorts',
    'medical compression tights',
    'medical compression tights women',
    'medical gear',
    'medical grade',
    'medical grade abdominal compression',
    'medical grade compression',
    'medical grade compression leggings',
    'medical grade compression leggings women',
    'medical grade compression pants',
    'medical grade compression tights',
    'medical grade fabric',
    'medical grade leggings',
    'medical jackets for women fitted',
    'medical leg compression',
    'medical pants for women',
    'medical shorts',
    'medical support',
    'medical ventilator',
    'medical ventilator machine',
    'medico compression',
    'medium compression pants',
    'medium size basketball',
    'medium support',
    'men capri leggings',
    'men capri tights',
    'men cold gear',
    'men cold gear compression',
    'men compression 3 4',
    'men compression 3 4 pants',
    'men compression leggings 3 4',
    'men compression nike',
    'men compression pants 3 4',
    'men compression pants blue',
    'men compression pants red',
    'men compression running pants',
    'men compression tights',
    'men compression tights 3 4',
    'men exercise tights',
    'men leggings 3 4',
    'men long shorts below the knee',
    'men nike compression pants',
    'men pants wicking fabric',
    'men red tights',
    'men running leggings warm',
    'men running pants cold weather',
    'men running thermal tights',
    'men s long rise travel pants',
    'men soccer wear',
    'men sport wear',
    'men thigh slimmer',
    'men tights legging',
    'men tights nike',
    'men track tights',
    'men winter tights nike',
    'men workout capri',
    'men workout tights',
    'mens 3/4 compression pants',
    'mens athletic leggings black',
    'mens athletic leggins',
    'mens ballet tights',
    'mens bicycle leg warmers',
    'mens black ballet tights',
    'mens brown leggings',
    'mens brown tights',
    'mens capri compression tights',
    'mens capri leggings for sports',
    'mens cold tights',
    'mens cold weather compression tights',
    'mens cold weather gear running',
    'mens compresion',
    'mens compression 3 4',
    'mens compression 3 4 leggings',
    'mens compression 3 4 pants',
    'mens compression 3 4 tights',
    'mens compression capri',
    'mens compression capri leggings',
    'mens compression capri pants',
    'mens compression capris',
    'mens compression cold gear',
    'mens compression leggings 3 4',
    'mens compression leggings capri',
    'mens compression leggings cold weather',
    'mens compression long',
    'mens compression pants 3 4',
    'mens compression pants 3/4',
    'mens compression pants 3/4 length',
    'mens compression pants capri',
    'mens compression pants cold',
    'mens compression pants running',
    'mens compression recovery pants',
    'mens compression running pants',
    'mens compression shorts 3 4',
    'mens compression shorts 6 inch',
    'mens compression stomach',
    'mens compression tights',
    'mens compression tights 3 4',
    'mens compression tights basketball',
    'mens compression tights black',
    'mens compression tights cold weather'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

